$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2 (including its cell styles/number formats) down into the
# two new rows, 3 and 4, so the new rows inherit the same formatting as the
# existing data row.
$ws.Range("A2:G2").Copy($ws.Range("A3:G3")) | Out-Null
$ws.Range("A2:G2").Copy($ws.Range("A4:G4")) | Out-Null

# Row 3: new student "Diego Encina Poblete" / rut "21.092.939-8"
$ws.Range("E3").Value = "Diego Encina Poblete"
$ws.Range("D3").Value = "21.092.939-8"

# Existing row 2's rut was stored as a plain number; re-enter it as the
# properly formatted rut text value (with punctuation), same as the rest.
$ws.Range("D2").Value = "21.510.487-7"

# Row 4: new student "Felipe Ignacio Tapia Diaz" / rut "21.126.460-8"
$ws.Range("D4").Value = "21.126.460-8"
$ws.Range("E4").Value = "Felipe Ignacio Tapia Diaz"

# Move the active selection to the last filled-in cell.
$ws.Range("G4").Select() | Out-Null
